$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 9, shifting existing rows 9-31 down to 10-32.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new weekly data record.
$ws.Cells.Item(9, 1).Value = 4
$ws.Cells.Item(9, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(9, 3).Value = "Los Lagos"
$ws.Cells.Item(9, 4).Value = 44838
$ws.Cells.Item(9, 5).Value = 10
$ws.Cells.Item(9, 6).Value = 100112013
$ws.Cells.Item(9, 7).Value = "Alcachofa"
$ws.Cells.Item(9, 8).Value = "Española"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 150
$ws.Cells.Item(9, 11).Value = 12000
$ws.Cells.Item(9, 12).Value = 12000
$ws.Cells.Item(9, 13).Value = 12000
$ws.Cells.Item(9, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(9, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(9, 16).Value = 400
$ws.Cells.Item(9, 17).Value = 30
$ws.Cells.Item(9, 18).Value = "Hortaliza"
